$wb = $excel.ActiveWorkbook

# Sheet: 展览 (sheet1) - column F ("想去人数") updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(2, 6).Value = 39
$ws1.Cells.Item(4, 6).Value = 4857
$ws1.Cells.Item(5, 6).Value = 0
$ws1.Cells.Item(6, 6).Value = 164
$ws1.Cells.Item(8, 6).Value = 0
$ws1.Cells.Item(9, 6).Value = 99
$ws1.Cells.Item(12, 6).Value = 1219
$ws1.Cells.Item(13, 6).Value = 123
$ws1.Cells.Item(14, 6).Value = 259
$ws1.Cells.Item(15, 6).Value = 201
$ws1.Cells.Item(17, 6).Value = 1
$ws1.Cells.Item(19, 6).Value = 116
$ws1.Cells.Item(20, 6).Value = 4125
$ws1.Cells.Item(21, 6).Value = 6433
$ws1.Cells.Item(22, 6).Value = 39
$ws1.Cells.Item(24, 6).Value = 0
$ws1.Cells.Item(25, 6).Value = 545
$ws1.Cells.Item(27, 6).Value = 0
$ws1.Cells.Item(28, 6).Value = 0
$ws1.Cells.Item(34, 6).Value = 0
$ws1.Cells.Item(35, 6).Value = 310
$ws1.Cells.Item(36, 6).Value = 324
$ws1.Cells.Item(38, 6).Value = 187
$ws1.Cells.Item(39, 6).Value = 14
$ws1.Cells.Item(40, 6).Value = 1578
$ws1.Cells.Item(41, 6).Value = 980
$ws1.Cells.Item(42, 6).Value = 49
$ws1.Cells.Item(43, 6).Value = 0
$ws1.Cells.Item(45, 6).Value = 0
$ws1.Cells.Item(46, 6).Value = 485
$ws1.Cells.Item(47, 6).Value = 0
$ws1.Cells.Item(48, 6).Value = 0
$ws1.Cells.Item(49, 6).Value = 596

# Sheet: 演出 (sheet2) - column F updates
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(2, 6).Value = 0

# Sheet: 全部类型 (sheet4) - column F updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(2, 6).Value = 39
$ws4.Cells.Item(3, 6).Value = 240
$ws4.Cells.Item(4, 6).Value = 4857
$ws4.Cells.Item(5, 6).Value = 212
$ws4.Cells.Item(7, 6).Value = 126
$ws4.Cells.Item(10, 6).Value = 99
$ws4.Cells.Item(11, 6).Value = 769
$ws4.Cells.Item(19, 6).Value = 116
$ws4.Cells.Item(20, 6).Value = 4125
$ws4.Cells.Item(21, 6).Value = 6433
$ws4.Cells.Item(22, 6).Value = 0
$ws4.Cells.Item(24, 6).Value = 0
$ws4.Cells.Item(25, 6).Value = 0
$ws4.Cells.Item(27, 6).Value = 4009
$ws4.Cells.Item(28, 6).Value = 412
$ws4.Cells.Item(29, 6).Value = 0
$ws4.Cells.Item(31, 6).Value = 2607
$ws4.Cells.Item(34, 6).Value = 153
$ws4.Cells.Item(35, 6).Value = 310
$ws4.Cells.Item(36, 6).Value = 324
$ws4.Cells.Item(37, 6).Value = 0
$ws4.Cells.Item(38, 6).Value = 187
$ws4.Cells.Item(39, 6).Value = 14
$ws4.Cells.Item(40, 6).Value = 1578
$ws4.Cells.Item(41, 6).Value = 0
$ws4.Cells.Item(42, 6).Value = 0
$ws4.Cells.Item(43, 6).Value = 0
$ws4.Cells.Item(44, 6).Value = 0
$ws4.Cells.Item(45, 6).Value = 504
$ws4.Cells.Item(46, 6).Value = 0
$ws4.Cells.Item(47, 6).Value = 6
$ws4.Cells.Item(48, 6).Value = 80
$ws4.Cells.Item(49, 6).Value = 0
